$wb = $excel.ActiveWorkbook

# The same F-column updates apply to both the "展览" and "全部类型" sheets,
# which contain duplicated data.
$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F ("想去人数")
$updates = @{
    4  = 66
    6  = 7027
    7  = 200
    8  = 162
    9  = 1052
    10 = 424
    12 = 189
    13 = 602
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
